$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (row 1)
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Update the one-hot indicator values that shifted between rows/columns
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1

$ws.Range("C3").Value = 0
$ws.Range("F3").Value = 1

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1

$ws.Range("A7").Value = 1
$ws.Range("F7").Value = 0
